# Updated cryptos list on Mon Jun 12 08:42:08 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Row=2; D="25.893.01"; E="  +0.46%  "},
    @{Row=3; D="1.743.43"; E="  -0.56%  "},
    @{Row=4; D="0.9983"; E="  -0.32%  "},
    @{Row=5; D="227.06"; E="  -4.26%  "},
    @{Row=6; D="0.9993"; E="  -0.18%  "},
    @{Row=7; D="0.5144"; E="  +1.75%  "},
    @{Row=8; D="0.2807"; E="  +6.99%  "},
    @{Row=9; D="39.16"; E="  -3.55%  "},
    @{Row=10; D="0.06094"; E="  -1.80%  "},
    @{Row=11; D="1.745.82"; E="  -0.16%  "},
    @{Row=12; D="0.06992"; E="  +0.67%  "},
    @{Row=13; D="15.26"; E="  -1.23%  "},
    @{Row=14; D="0.6350"; E="  +5.34%  "},
    @{Row=15; D="4.495"; E="  +0.98%  "},
    @{Row=16; D="76.54"; E="  -2.53%  "},
    @{Row=17; D="0.9986"; E="  -0.27%  "},
    @{Row=18; D="0.9991"; E="  -0.20%  "},
    @{Row=19; D="25.891.44"; E="  +0.31%  "},
    @{Row=20; D="11.47"; E="  -1.51%  "},
    @{Row=21; D="0.000006580"; E="  -2.83%  "},
    @{Row=22; D="1.966.84"; E="  -0.42%  "},
    @{Row=23; D="4.096"; E="  +1.10%  "},
    @{Row=24; D="8.449"; E="  +3.05%  "},
    @{Row=25; D="5.126"; E="  -0.87%  "},
    @{Row=26; D="138.84"; E="  +0.67%  "},
    @{Row=27; D="1.509"},
    @{Row=28; D="1.831"; E="  +1.95%  "},
    @{Row=29; D="15.04"; E="  -0.31%  "},
    @{Row=30; D="102.98"; E="  +0.67%  "},
    @{Row=31; D="0.08307"; E="  +0.36%  "},
    @{Row=32; D="3.619"; E="  -2.12%  "},
    @{Row=33; D="3.412"; E="  +0.71%  "},
    @{Row=34; D="0.04382"; E="  +0.45%  "},
    @{Row=35; D="2.618"; E="  -1.04%  "},
    @{Row=36; D="0.9693"; E="  -3.13%  "},
    @{Row=37; D="0.6068"; E="  +1.03%  "},
    @{Row=38; D="2.665"; E="  -1.23%  "},
    @{Row=39; D="0.01559"; E="  +0.87%  "},
    @{Row=40; D="1.917"; E="  -2.07%  "},
    @{Row=41; D="0.9988"; E="  -0.23%  "},
    @{Row=42; D="100.62"; E="  -2.62%  "},
    @{Row=43; D="0.3843"; E="  +1.09%  "},
    @{Row=44; D="0.7219"; E="  -3.33%  "},
    @{Row=45; D="4.924"; E="  +0.97%  "},
    @{Row=46; D="0.05443"; E="  -0.69%  "},
    @{Row=47; D="6.322"; E="  +6.39%  "},
    @{Row=48; D="0.1105"; E="  +2.69%  "},
    @{Row=49; D="52.51"; E="  +1.20%  "},
    @{Row=50; D="29.74"; E="  -1.40%  "},
    @{Row=51; D="7.472"; E="  +0.65%  "}
)

foreach ($u in $updates) {
    $row = $u.Row
    $dCell = $ws.Cells.Item($row, 4)
    $dCell.NumberFormat = "@"
    $dCell.Value = $u.D

    if ($u.ContainsKey("E")) {
        $eCell = $ws.Cells.Item($row, 5)
        $eCell.NumberFormat = "@"
        $eCell.Value = $u.E
    }
}
